# Switzerland Challenge League - odds base update (02-05-2024 20:28)
# This script reorders several pairs of fixture rows that share an identical
# kickoff timestamp (the upstream feed flipped their order on refresh) and
# refreshes the odds for the three still-to-be-played fixtures that were
# re-ordered among themselves (rows 164-166), including a small odds tweak
# on the still-fixed row 163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rangeA = "B$r1`:AB$r1"
    $rangeB = "B$r2`:AB$r2"
    $valsA = $ws.Range($rangeA).Value2
    $valsB = $ws.Range($rangeB).Value2
    $ws.Range($rangeA).Value = $valsB
    $ws.Range($rangeB).Value = $valsA
}

# Pairs of rows whose full B:AB content (match id, teams, scores, odds, ...)
# swap places, while the sequential id in column A stays tied to the row.
Swap-Rows 2 3
Swap-Rows 13 14
Swap-Rows 28 29
Swap-Rows 39 40
Swap-Rows 47 48
Swap-Rows 78 79
Swap-Rows 90 91

# Rows 163-166 (all upcoming fixtures, no result yet) — fixtures 164/165/166
# got reshuffled and all four rows received refreshed odds.

# Row 163 stays the same fixture (AC Bellinzona vs Neuchatel Xamax) but with
# updated current odds.
$ws.Range("M163").Value = 2.8
$ws.Range("O163").Value = 2.55
$ws.Range("P163").Value = 0
$ws.Range("Q163").Value = 2.025
$ws.Range("R163").Value = 1.775

# Row 164 becomes the FC Baden vs FC Thun fixture (previously row 165)
$ws.Range("B164").Value = "'7617816"
$ws.Range("E164").Value = "FC Baden"
$ws.Range("F164").Value = "FC Thun"
$ws.Range("J164").Value = 5.5
$ws.Range("K164").Value = 4.5
$ws.Range("L164").Value = 1.5
$ws.Range("M164").Value = 6.5
$ws.Range("N164").Value = 4.75
$ws.Range("O164").Value = 1.45
$ws.Range("P164").Value = 1.25
$ws.Range("Q164").Value = 1.85
$ws.Range("R164").Value = 1.95
$ws.Range("S164").Value = 3.25
$ws.Range("T164").Value = 2.025
$ws.Range("U164").Value = 1.775

# Row 165 becomes the FC Vaduz vs Aarau fixture (previously row 166)
$ws.Range("B165").Value = "'7617813"
$ws.Range("E165").Value = "FC Vaduz"
$ws.Range("F165").Value = "Aarau"
$ws.Range("J165").Value = 2.1
$ws.Range("K165").Value = 3.75
$ws.Range("L165").Value = 3.1
$ws.Range("M165").Value = 2
$ws.Range("N165").Value = 3.6
$ws.Range("O165").Value = 3.5
$ws.Range("P165").Value = -0.5
$ws.Range("Q165").Value = 2
$ws.Range("R165").Value = 1.8
$ws.Range("T165").Value = 1.85
$ws.Range("U165").Value = 1.95

# Row 166 becomes the FC Sion vs Wil 1900 fixture (previously row 164)
$ws.Range("B166").Value = "'7617814"
$ws.Range("E166").Value = "FC Sion"
$ws.Range("F166").Value = "Wil 1900"
$ws.Range("J166").Value = 1.6
$ws.Range("K166").Value = 4
$ws.Range("L166").Value = 5.25
$ws.Range("M166").Value = 1.533
$ws.Range("N166").Value = 4.333
$ws.Range("O166").Value = 6
$ws.Range("P166").Value = -1
$ws.Range("Q166").Value = 1.85
$ws.Range("R166").Value = 1.95
$ws.Range("S166").Value = 2.75
$ws.Range("T166").Value = 1.8
$ws.Range("U166").Value = 2
